$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 8
$ws.Range("A8").Value = 42061
$ws.Range("A8").NumberFormat = "yyyy/mm/dd;@"
$ws.Range("C8").Value = "Requirement Specifier"
$ws.Range("D8").Value = "Activitets Diagram og Domænemodel for UC-5"
$ws.Range("B8").Value = "3 hrs?"

# Add new row 9
$ws.Range("A9").Value = 42062
$ws.Range("A9").NumberFormat = "yyyy/mm/dd;@"
$ws.Range("B9").Value = "4,30 hrs?"

# Correct the existing row (row 5, "System Sekvens Diagram") time value
$ws.Range("B5").Value = "2,30 hrs?"

$ws.Range("C9").Value = "Test analyst"
$ws.Range("D9").Value = "Udarbejdelse af test til OC-2"

# Update the selection to match the new active cell
$ws.Range("D10").Select()
